$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(0.5553305257193456, 0.6136273584998631, 0.5228536874213209, 0.723086224057215, 0.4805816767817249, 14)
    3  = @(0.3621708991783025, 0.4983652827815714, 0.342784039812318, 0.585477616832888, 0.4788015972899659, 13)
    4  = @(0.322109523816475, 0.4916759303336116, 0.3220049964818212, 0.5674548409184833, 0.4879461048650333, 12)
    5  = @(0.4479694299613027, 0.5810026869938167, 0.4406117130697931, 0.6637858939972987, 0.5137398302438039, 11)
    6  = @(0.4056295499570255, 0.5539856660635349, 0.4125774920458339, 0.6423219535761128, 0.5249784760111545, 10)
    7  = @(0.3136617854706863, 0.4838333121955758, 0.3116578194982599, 0.5582632170385757, 0.4898299366237236, 9)
    8  = @(0.3613000660075406, 0.5407655076793962, 0.3759763669136598, 0.6131691177103261, 0.5296237254251971, 8)
    9  = @(0.5028015250901602, 0.5533246075941491, 0.4051903875053499, 0.6365456680438175, 0.4216371064289196, 7)
    10 = @(0.4053476105038098, 0.459336539820115, 0.2581058573136952, 0.5080411964729782, 0.33549814659258, 6)
    11 = @(0.4015449179333902, 0.4597420797009394, 0.2670924375257908, 0.516809865933102, 0.3637549250657264, 5)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Range("B$row").Value = $rowValues[0]
    $ws.Range("C$row").Value = $rowValues[1]
    $ws.Range("D$row").Value = $rowValues[2]
    $ws.Range("E$row").Value = $rowValues[3]
    $ws.Range("F$row").Value = $rowValues[4]
    $ws.Range("G$row").Value = $rowValues[5]
}
